# Daily auto push: 2026-01-20 02:29 UTC
#
# A new observation row for 2026/01/20 (weekday 火, hour 9, value 15) was
# recorded. It belongs right after the existing "2026/01/20" rows (which
# currently end at row 660) and before the "2026/12/29" rows, so a single
# row is inserted at row 661, pushing the old rows 661-702 down to 662-703.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert one blank row at position 661 (existing rows 661-702 shift down to
# 662-703, and the sheet's used range / dimension grows from D702 to D703).
$ws.Rows.Item(661).Insert()

# Column A holds dates stored as literal text (e.g. "2026/01/20"), not as
# real Excel date serials. Prefixing with an apostrophe forces the value to
# stay text instead of being auto-converted to a date; resetting the style
# back to "Normal" afterwards clears the quote-prefix formatting flag so the
# cell's style matches the rest of the (unstyled) data rows.
$ws.Range("A661").Value = "'2026/01/20"
$ws.Range("A661").Style = "Normal"

$ws.Range("B661").Value = "火"
$ws.Range("C661").Value = 9
$ws.Range("D661").Value = 15
